$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Participant query" cell (B2, row for the ParticipantsTab) is replaced
# with a more elaborate Cypher query that goes through diagnosis / file /
# genomic_info and returns a sorted, de-duplicated sample-id list.
$newQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE f.file_type in ['HTML']
with p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id LIMIT 100
'@

$ws.Range("B2").Value = $newQuery

# The longer text needs a taller row (12 lines -> 18 lines).
$ws.Rows.Item(2).RowHeight = 279

# View state: scrolled down so row 3 is at the top, with B4 selected.
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("B4").Select()
